# Corrects the IFRS financial summary rows (2-9) on the "company_list" sheet
# to the restated figures, replacing erroneous duplicated/merged totals and
# removing columns that no longer apply for the forecast rows (7-9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> corrected value (applied with Cells.Item(row, col) to avoid locale/
# culture issues with A1-style Range addressing of two-letter columns)
$updates = @(
    @{Row=2; Col=4; Value=79549}  # D2
    @{Row=2; Col=5; Value=1584}  # E2
    @{Row=2; Col=6; Value=1584}  # F2
    @{Row=2; Col=7; Value=1554}  # G2
    @{Row=2; Col=8; Value=1175}  # H2
    @{Row=2; Col=9; Value=1175}  # I2
    @{Row=2; Col=10; Value=0}  # J2
    @{Row=2; Col=11; Value=84639}  # K2
    @{Row=2; Col=12; Value=66249}  # L2
    @{Row=2; Col=13; Value=18390}  # M2
    @{Row=2; Col=14; Value=18390}  # N2
    @{Row=2; Col=16; Value=602}  # P2
    @{Row=2; Col=17; Value=44}  # Q2
    @{Row=2; Col=18; Value=-2280}  # R2
    @{Row=2; Col=19; Value=1921}  # S2
    @{Row=2; Col=20; Value=46}  # T2
    @{Row=2; Col=22; Value=0}  # V2
    @{Row=2; Col=23; Value=1.99}  # W2
    @{Row=2; Col=24; Value=1.48}  # X2
    @{Row=2; Col=25; Value=7.14}  # Y2
    @{Row=2; Col=26; Value=1.46}  # Z2
    @{Row=2; Col=27; Value=360.25}  # AA2
    @{Row=2; Col=28; Value=3014.21}  # AB2
    @{Row=2; Col=29; Value=976}  # AC2
    @{Row=2; Col=30; Value=11.07}  # AD2
    @{Row=2; Col=31; Value=16009}  # AE2
    @{Row=2; Col=32; Value=0.67}  # AF2
    @{Row=2; Col=33; Value=225}  # AG2
    @{Row=2; Col=34; Value=2.08}  # AH2
    @{Row=2; Col=35; Value=22}  # AI2
    @{Row=2; Col=36; Value=120369116}  # AJ2
    @{Row=3; Col=4; Value=84301}  # D3
    @{Row=3; Col=5; Value=2453}  # E3
    @{Row=3; Col=6; Value=2453}  # F3
    @{Row=3; Col=7; Value=2416}  # G3
    @{Row=3; Col=8; Value=1865}  # H3
    @{Row=3; Col=9; Value=1865}  # I3
    @{Row=3; Col=11; Value=89785}  # K3
    @{Row=3; Col=12; Value=69617}  # L3
    @{Row=3; Col=13; Value=20168}  # M3
    @{Row=3; Col=14; Value=20168}  # N3
    @{Row=3; Col=16; Value=602}  # P3
    @{Row=3; Col=17; Value=1142}  # Q3
    @{Row=3; Col=18; Value=-1434}  # R3
    @{Row=3; Col=19; Value=-366}  # S3
    @{Row=3; Col=20; Value=3}  # T3
    @{Row=3; Col=22; Value=0}  # V3
    @{Row=3; Col=23; Value=2.91}  # W3
    @{Row=3; Col=24; Value=2.21}  # X3
    @{Row=3; Col=25; Value=9.67}  # Y3
    @{Row=3; Col=26; Value=2.14}  # Z3
    @{Row=3; Col=27; Value=345.19}  # AA3
    @{Row=3; Col=28; Value=3309.63}  # AB3
    @{Row=3; Col=29; Value=1549}  # AC3
    @{Row=3; Col=30; Value=9.130000000000001}  # AD3
    @{Row=3; Col=31; Value=17557}  # AE3
    @{Row=3; Col=32; Value=0.8100000000000001}  # AF3
    @{Row=3; Col=33; Value=350}  # AG3
    @{Row=3; Col=34; Value=2.47}  # AH3
    @{Row=3; Col=35; Value=21.56}  # AI3
    @{Row=3; Col=36; Value=120369116}  # AJ3
    @{Row=4; Col=4; Value=85597}  # D4
    @{Row=4; Col=5; Value=2072}  # E4
    @{Row=4; Col=6; Value=2072}  # F4
    @{Row=4; Col=7; Value=2083}  # G4
    @{Row=4; Col=8; Value=1600}  # H4
    @{Row=4; Col=9; Value=1600}  # I4
    @{Row=4; Col=11; Value=95811}  # K4
    @{Row=4; Col=12; Value=74695}  # L4
    @{Row=4; Col=13; Value=21116}  # M4
    @{Row=4; Col=14; Value=21116}  # N4
    @{Row=4; Col=16; Value=602}  # P4
    @{Row=4; Col=17; Value=4171}  # Q4
    @{Row=4; Col=18; Value=-3692}  # R4
    @{Row=4; Col=19; Value=-504}  # S4
    @{Row=4; Col=20; Value=27}  # T4
    @{Row=4; Col=22; Value=0}  # V4
    @{Row=4; Col=23; Value=2.42}  # W4
    @{Row=4; Col=24; Value=1.87}  # X4
    @{Row=4; Col=25; Value=7.75}  # Y4
    @{Row=4; Col=26; Value=1.73}  # Z4
    @{Row=4; Col=27; Value=353.74}  # AA4
    @{Row=4; Col=28; Value=3467.19}  # AB4
    @{Row=4; Col=29; Value=1330}  # AC4
    @{Row=4; Col=30; Value=8.57}  # AD4
    @{Row=4; Col=31; Value=18383}  # AE4
    @{Row=4; Col=32; Value=0.62}  # AF4
    @{Row=4; Col=33; Value=325}  # AG4
    @{Row=4; Col=34; Value=2.85}  # AH4
    @{Row=4; Col=35; Value=23.33}  # AI4
    @{Row=4; Col=36; Value=120369116}  # AJ4
    @{Row=5; Col=4; Value=91168}  # D5
    @{Row=5; Col=5; Value=1756}  # E5
    @{Row=5; Col=6; Value=1756}  # F5
    @{Row=5; Col=7; Value=1734}  # G5
    @{Row=5; Col=8; Value=1330}  # H5
    @{Row=5; Col=9; Value=1330}  # I5
    @{Row=5; Col=11; Value=100653}  # K5
    @{Row=5; Col=12; Value=79023}  # L5
    @{Row=5; Col=13; Value=21630}  # M5
    @{Row=5; Col=14; Value=21630}  # N5
    @{Row=5; Col=16; Value=602}  # P5
    @{Row=5; Col=17; Value=5596}  # Q5
    @{Row=5; Col=18; Value=-5790}  # R5
    @{Row=5; Col=19; Value=-475}  # S5
    @{Row=5; Col=20; Value=9}  # T5
    @{Row=5; Col=22; Value=0}  # V5
    @{Row=5; Col=23; Value=1.93}  # W5
    @{Row=5; Col=24; Value=1.46}  # X5
    @{Row=5; Col=25; Value=6.22}  # Y5
    @{Row=5; Col=26; Value=1.35}  # Z5
    @{Row=5; Col=27; Value=365.34}  # AA5
    @{Row=5; Col=28; Value=3552.64}  # AB5
    @{Row=5; Col=29; Value=1105}  # AC5
    @{Row=5; Col=30; Value=9.949999999999999}  # AD5
    @{Row=5; Col=31; Value=18830}  # AE5
    @{Row=5; Col=32; Value=0.58}  # AF5
    @{Row=5; Col=33; Value=300}  # AG5
    @{Row=5; Col=34; Value=2.73}  # AH5
    @{Row=5; Col=35; Value=25.91}  # AI5
    @{Row=5; Col=36; Value=120369116}  # AJ5
    @{Row=6; Col=4; Value=98925}  # D6
    @{Row=6; Col=5; Value=1436}  # E6
    @{Row=6; Col=6; Value=1436}  # F6
    @{Row=6; Col=7; Value=1344}  # G6
    @{Row=6; Col=8; Value=1029}  # H6
    @{Row=6; Col=9; Value=1029}  # I6
    @{Row=6; Col=11; Value=107509}  # K6
    @{Row=6; Col=12; Value=85132}  # L6
    @{Row=6; Col=13; Value=22376}  # M6
    @{Row=6; Col=14; Value=22376}  # N6
    @{Row=6; Col=16; Value=602}  # P6
    @{Row=6; Col=17; Value=3091}  # Q6
    @{Row=6; Col=18; Value=-2194}  # R6
    @{Row=6; Col=19; Value=-445}  # S6
    @{Row=6; Col=20; Value=10}  # T6
    @{Row=6; Col=22; Value=0}  # V6
    @{Row=6; Col=23; Value=1.45}  # W6
    @{Row=6; Col=24; Value=1.04}  # X6
    @{Row=6; Col=25; Value=4.68}  # Y6
    @{Row=6; Col=26; Value=0.99}  # Z6
    @{Row=6; Col=27; Value=380.46}  # AA6
    @{Row=6; Col=28; Value=3676.59}  # AB6
    @{Row=6; Col=29; Value=855}  # AC6
    @{Row=6; Col=30; Value=10.13}  # AD6
    @{Row=6; Col=31; Value=19480}  # AE6
    @{Row=6; Col=32; Value=0.44}  # AF6
    @{Row=6; Col=33; Value=275}  # AG6
    @{Row=6; Col=34; Value=3.18}  # AH6
    @{Row=6; Col=35; Value=30.69}  # AI6
    @{Row=6; Col=36; Value=120369116}  # AJ6
    @{Row=7; Col=33; Value=475}  # AG7
    @{Row=7; Col=34; Value=5.5}  # AH7
    @{Row=8; Col=33; Value=490}  # AG8
    @{Row=8; Col=34; Value=5.68}  # AH8
    @{Row=9; Col=33; Value=525}  # AG9
    @{Row=9; Col=34; Value=6.08}  # AH9
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# Cells that are blank in the corrected data (the source columns were dropped
# for these rows) - remove them entirely rather than leaving a 0/empty value
$clears = @(
    @{Row=2; Col=15}  # O2
    @{Row=2; Col=21}  # U2
    @{Row=3; Col=10}  # J3
    @{Row=3; Col=15}  # O3
    @{Row=3; Col=21}  # U3
    @{Row=4; Col=10}  # J4
    @{Row=4; Col=15}  # O4
    @{Row=4; Col=21}  # U4
    @{Row=5; Col=10}  # J5
    @{Row=5; Col=15}  # O5
    @{Row=5; Col=21}  # U5
    @{Row=6; Col=21}  # U6
    @{Row=7; Col=4}  # D7
    @{Row=7; Col=5}  # E7
    @{Row=7; Col=7}  # G7
    @{Row=7; Col=8}  # H7
    @{Row=7; Col=9}  # I7
    @{Row=7; Col=11}  # K7
    @{Row=7; Col=12}  # L7
    @{Row=7; Col=13}  # M7
    @{Row=7; Col=14}  # N7
    @{Row=7; Col=16}  # P7
    @{Row=7; Col=17}  # Q7
    @{Row=7; Col=18}  # R7
    @{Row=7; Col=19}  # S7
    @{Row=7; Col=20}  # T7
    @{Row=7; Col=21}  # U7
    @{Row=7; Col=23}  # W7
    @{Row=7; Col=24}  # X7
    @{Row=7; Col=25}  # Y7
    @{Row=7; Col=26}  # Z7
    @{Row=7; Col=27}  # AA7
    @{Row=7; Col=29}  # AC7
    @{Row=7; Col=30}  # AD7
    @{Row=7; Col=31}  # AE7
    @{Row=7; Col=32}  # AF7
    @{Row=7; Col=35}  # AI7
    @{Row=8; Col=4}  # D8
    @{Row=8; Col=5}  # E8
    @{Row=8; Col=7}  # G8
    @{Row=8; Col=8}  # H8
    @{Row=8; Col=9}  # I8
    @{Row=8; Col=11}  # K8
    @{Row=8; Col=12}  # L8
    @{Row=8; Col=13}  # M8
    @{Row=8; Col=14}  # N8
    @{Row=8; Col=16}  # P8
    @{Row=8; Col=17}  # Q8
    @{Row=8; Col=18}  # R8
    @{Row=8; Col=19}  # S8
    @{Row=8; Col=20}  # T8
    @{Row=8; Col=21}  # U8
    @{Row=8; Col=23}  # W8
    @{Row=8; Col=24}  # X8
    @{Row=8; Col=25}  # Y8
    @{Row=8; Col=26}  # Z8
    @{Row=8; Col=27}  # AA8
    @{Row=8; Col=29}  # AC8
    @{Row=8; Col=30}  # AD8
    @{Row=8; Col=31}  # AE8
    @{Row=8; Col=32}  # AF8
    @{Row=8; Col=35}  # AI8
    @{Row=9; Col=4}  # D9
    @{Row=9; Col=5}  # E9
    @{Row=9; Col=7}  # G9
    @{Row=9; Col=8}  # H9
    @{Row=9; Col=9}  # I9
    @{Row=9; Col=11}  # K9
    @{Row=9; Col=12}  # L9
    @{Row=9; Col=13}  # M9
    @{Row=9; Col=14}  # N9
    @{Row=9; Col=16}  # P9
    @{Row=9; Col=17}  # Q9
    @{Row=9; Col=18}  # R9
    @{Row=9; Col=19}  # S9
    @{Row=9; Col=20}  # T9
    @{Row=9; Col=21}  # U9
    @{Row=9; Col=23}  # W9
    @{Row=9; Col=24}  # X9
    @{Row=9; Col=25}  # Y9
    @{Row=9; Col=26}  # Z9
    @{Row=9; Col=27}  # AA9
    @{Row=9; Col=29}  # AC9
    @{Row=9; Col=30}  # AD9
    @{Row=9; Col=31}  # AE9
    @{Row=9; Col=32}  # AF9
    @{Row=9; Col=35}  # AI9
)

foreach ($c in $clears) {
    $ws.Cells.Item($c.Row, $c.Col).ClearContents()
}

